$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.01848533333333334
$ws.Range("N2").Value = 0.05545600000000001
$ws.Range("O2").Value = 0.001625201930372746
$ws.Range("P2").Value = 0.001625201930372746
$ws.Range("Q2").Value = 0.0006341085511111112
$ws.Range("R2").Value = 0.00570697696
$ws.Range("S2").Value = 0.001625201930372746
$ws.Range("T2").Value = 0.001625201930372746

# Row 3 updates
$ws.Range("O3").Value = 0.002698334581238102
$ws.Range("P3").Value = 0.002698334581238102
$ws.Range("S3").Value = 0.002698334581238102
$ws.Range("T3").Value = 0.002698334581238102

# Row 4 updates
$ws.Range("M4").Value = 11.32499966666667
$ws.Range("N4").Value = 33.974999
$ws.Range("O4").Value = 0.9956764634883892
$ws.Range("P4").Value = 0.995676463488389
$ws.Range("Q4").Value = 0.3884852385655556
$ws.Range("R4").Value = 3.49636714709
$ws.Range("S4").Value = 0.9956764634883892
$ws.Range("T4").Value = 0.995676463488389
